$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 505 and 506, pushing existing rows 505-541 down to 507-543.
$ws.Rows("505:506").Insert()

# Constant column values shared by every data row in this sheet.
$constA = 4
$constB = "Feria Lagunitas de Puerto Montt"
$constC = "Los Lagos"
$constE = 10
$constF = "Fruta"
$constG = 100102
$constH = "Cítricos"
$constI = 100102006
$constJ = "Pomelo"
$constQ = "`$/caja 14 kilos empedrada"
$constR = "Región de O'Higgins"
$constT = 14

function Set-DataRow($RowNum, $Fecha, $Variedad, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $PrecioKg) {
    $ws.Cells.Item($RowNum, 1).Value2 = $constA
    $ws.Cells.Item($RowNum, 2).Value2 = $constB
    $ws.Cells.Item($RowNum, 3).Value2 = $constC
    $ws.Cells.Item($RowNum, 4).Value2 = $Fecha
    $ws.Cells.Item($RowNum, 5).Value2 = $constE
    $ws.Cells.Item($RowNum, 6).Value2 = $constF
    $ws.Cells.Item($RowNum, 7).Value2 = $constG
    $ws.Cells.Item($RowNum, 8).Value2 = $constH
    $ws.Cells.Item($RowNum, 9).Value2 = $constI
    $ws.Cells.Item($RowNum, 10).Value2 = $constJ
    $ws.Cells.Item($RowNum, 11).Value2 = $Variedad
    $ws.Cells.Item($RowNum, 12).Value2 = $Calidad
    $ws.Cells.Item($RowNum, 13).Value2 = $Volumen
    $ws.Cells.Item($RowNum, 14).Value2 = $PrecioMin
    $ws.Cells.Item($RowNum, 15).Value2 = $PrecioMax
    $ws.Cells.Item($RowNum, 16).Value2 = $PrecioProm
    $ws.Cells.Item($RowNum, 17).Value2 = $constQ
    $ws.Cells.Item($RowNum, 18).Value2 = $constR
    $ws.Cells.Item($RowNum, 19).Value2 = $PrecioKg
    $ws.Cells.Item($RowNum, 20).Value2 = $constT
}

# New row 505: Start Ruby / Primera
Set-DataRow 505 45013 "Start Ruby" "Primera" 120 14000 15000 14500 1036

# New row 506: Start Ruby / Segunda
Set-DataRow 506 45013 "Start Ruby" "Segunda" 60 12000 12000 12000 857

# Apply the date number format (matches style used by the other "Fecha" column cells)
$dateFmt = $ws.Range("D504").NumberFormat
$ws.Range("D505:D506").NumberFormat = $dateFmt
